$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.3082271070678201,
    -2.241727458159763,
    0.1152157549956138,
    0.02297854259031917,
    0.1796463729455963,
    0.4106662722775318,
    1.924429454498304,
    1.151440957050255,
    0.307601050108237,
    0.7295210035792458,
    0.2690475247287556,
    0.6408324837877148,
    -1.07531867879654,
    0.6681140493286761,
    33.77994876784079,
    53.281961965732
)

for ($row = 2; $row -le 26; $row++) {
    for ($col = 2; $col -le 17; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 2]
    }
}
